$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 1.29
$ws.Range("H3").Value = 4.85
$ws.Range("K3").Value = 2.37
$ws.Range("N3").Value = 8
$ws.Range("O3").Value = 1.23
$ws.Range("P3").Value = 3.35
$ws.Range("Q3").Value = 1.7
$ws.Range("U3").Value = 2.18
$ws.Range("V3").Value = 1.53
$ws.Range("W3").Value = 6.2
$ws.Range("X3").Value = 5.6
$ws.Range("AA3").Value = 11.75
$ws.Range("AB3").Value = 35
$ws.Range("AC3").Value = 11.5
$ws.Range("AD3").Value = 10
$ws.Range("AF3").Value = 150
$ws.Range("AG3").Value = 21
$ws.Range("AI3").Value = 32
$ws.Range("AN3").Value = 2.95
$ws.Range("AO3").Value = 5.5
$ws.Range("AP3").Value = 18
$ws.Range("AQ3").Value = 15
$ws.Range("AR3").Value = 50
$ws.Range("AS3").Value = 2.87
$ws.Range("AU3").Value = 110
$ws.Range("BB3").Value = 300
$ws.Range("J6").Value = 2.9
$ws.Range("O6").Value = 1.5
$ws.Range("P6").Value = 2.42
$ws.Range("U6").Value = 2.05
$ws.Range("AE6").Value = 17
$ws.Range("AF6").Value = 110
$ws.Range("AG6").Value = 8
$ws.Range("AH6").Value = 17
$ws.Range("AO6").Value = 12
$ws.Range("AT6").Value = 7.6
$ws.Range("J7").Value = 2.32
$ws.Range("K7").Value = 2.1
$ws.Range("L7").Value = 4.75
$ws.Range("N7").Value = 9.050000000000001
$ws.Range("P7").Value = 2.7
$ws.Range("S7").Value = 1.42
$ws.Range("T7").Value = 2.47
$ws.Range("U7").Value = 1.98
$ws.Range("V7").Value = 1.65
$ws.Range("W7").Value = 5.9
$ws.Range("AA7").Value = 15.5
$ws.Range("AB7").Value = 35
$ws.Range("AE7").Value = 19
$ws.Range("AF7").Value = 110
$ws.Range("AJ7").Value = 70
$ws.Range("AP7").Value = 20
$ws.Range("AQ7").Value = 30
$ws.Range("AS7").Value = 2.45
$ws.Range("AT7").Value = 8
$ws.Range("AU7").Value = 90
$ws.Range("AV7").Value = 5.9
